# "Whiteboard drawings from lectures as separate files."
#
# The deck ended with four whiteboard-photo slides (SlideID 262, 261,
# 263, 267 -> ppt/slides/slide6.xml, slide7.xml, slide8.xml, slide9.xml).
# Keep only one of them (SlideID 261 / slide7.xml) and drop the rest
# (SlideID 262 / slide6.xml, SlideID 263 / slide8.xml, SlideID 267 /
# slide9.xml) now that those whiteboard drawings live as separate files
# instead of as extra slides in the deck.

$p = $ppt.ActivePresentation

$idsToRemove = @(262, 263, 267)

foreach ($targetId in $idsToRemove) {
    for ($i = $p.Slides.Count; $i -ge 1; $i--) {
        $slide = $p.Slides.Item($i)
        if ($slide.SlideID -eq $targetId) {
            $slide.Delete()
            break
        }
    }
}
